# Bugifixed QoQ Visualizations and a typo in the evaluation objects
# Updates the forecast-error statistics (ME, MAE, MSE, RMSE, SE) for
# quarters Q0-Q9 (rows 2-11, columns B-F) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2"  = -0.1250617157922688
    "C2"  = 0.5555183432352195
    "D2"  = 0.6995838223442629
    "E2"  = 0.8364112758351975
    "F2"  = 0.8354049424957229

    "B3"  = 0.3719515601421958
    "C3"  = 0.7459951107692967
    "D3"  = 1.522418059711977
    "E3"  = 1.233863063598217
    "F3"  = 1.188656758824574

    "B4"  = 0.4541585807807824
    "C4"  = 1.04916654926496
    "D4"  = 3.764402933008068
    "E4"  = 1.940206930460787
    "F4"  = 1.906265481234884

    "B5"  = 0.422918388690988
    "C5"  = 1.13962110842592
    "D5"  = 4.720541585393809
    "E5"  = 2.172680737106538
    "F5"  = 2.15416197883831

    "B6"  = 0.2396519184303125
    "C6"  = 0.962147624362754
    "D6"  = 4.10586913343047
    "E6"  = 2.026294434042217
    "F6"  = 2.034306107226543

    "B7"  = 0.2900791048022007
    "C7"  = 1.030905847179739
    "D7"  = 5.137395499217558
    "E7"  = 2.266582338945037
    "F7"  = 2.28174897574234

    "B8"  = 0.2209419343134404
    "C8"  = 1.046751655212453
    "D8"  = 5.254410289209334
    "E8"  = 2.292250049451266
    "F8"  = 2.316952695478674

    "B9"  = 0.1370776462960769
    "C9"  = 1.61474521251127
    "D9"  = 10.23233324787276
    "E9"  = 3.198801845671713
    "F9"  = 3.300673542053086

    "B10" = -0.7038664495325276
    "C10" = 1.239373893761973
    "D10" = 7.311748838265475
    "E10" = 2.704024563177168
    "F10" = 2.752033765031871

    "B11" = 0.2007429256103617
    "C11" = 0.6416278403412059
    "D11" = 0.4618474333756422
    "E11" = 0.6795935795574015
    "F11" = 0.7259043593968103
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
